# Updated cryptos list (prices in column D, 1h volume change % in column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.643.61"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.565.70"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.39"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.509"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.90"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.28%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.789.02"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.568.92"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.668.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.78"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.67"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.78"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.405.48"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.04"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.518"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.769"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.95"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.23"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.701.34"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.859"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.75"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.42%  "
$ws.Range("E51").Value = "  -0.45%  "
